$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,9

$data[0,0] = 3
$data[0,1] = 'control'
$data[0,2] = 'clim'
$data[0,3] = '[''ifs-fesom'', ''icon'']'
$data[0,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[0,5] = 'decadal/control_EERIE_clim.zarr'
$data[0,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[0,7] = ''
$data[0,8] = '[''period'', ''time_filter'']'
$data[1,0] = 8
$data[1,1] = 'control'
$data[1,2] = 'trend'
$data[1,3] = '[''ifs-fesom'', ''icon'']'
$data[1,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[1,5] = 'decadal/control_EERIE_trend.zarr'
$data[1,6] = '[''value (variable with hatching over variable_pvalue > 0.05)'']'
$data[1,7] = ''
$data[1,8] = '[''period'', ''time_filter'']'
$data[2,0] = 16
$data[2,1] = 'control'
$data[2,2] = 'ts'
$data[2,3] = '[''ifs-fesom'', ''icon'']'
$data[2,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[2,5] = 'time_series/control_EERIE_IPCC_ts.zarr'
$data[2,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[2,7] = 'IPCC'
$data[2,8] = '[''region'', ''time_filter'']'
$data[3,0] = 17
$data[3,1] = 'control'
$data[3,2] = 'ts'
$data[3,3] = '[''ifs-fesom'', ''icon'']'
$data[3,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[3,5] = 'time_series/control_EERIE_EDDY_ts.zarr'
$data[3,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[3,7] = 'EDDY'
$data[3,8] = '[''region'', ''time_filter'']'
$data[4,0] = 0
$data[4,1] = 'future'
$data[4,2] = 'clim'
$data[4,3] = '[''ifs-fesom'', ''icon'']'
$data[4,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[4,5] = 'decadal/future_EERIE_clim.zarr'
$data[4,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[4,7] = ''
$data[4,8] = '[''period'', ''time_filter'']'
$data[5,0] = 5
$data[5,1] = 'future'
$data[5,2] = 'trend'
$data[5,3] = '[''ifs-fesom'', ''icon'']'
$data[5,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[5,5] = 'decadal/future_EERIE_trend.zarr'
$data[5,6] = '[''value (variable with hatching over variable_pvalue > 0.05)'']'
$data[5,7] = ''
$data[5,8] = '[''period'', ''time_filter'']'
$data[6,0] = 10
$data[6,1] = 'future'
$data[6,2] = 'ts'
$data[6,3] = '[''ifs-fesom'', ''icon'']'
$data[6,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[6,5] = 'time_series/future_EERIE_IPCC_ts.zarr'
$data[6,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[6,7] = 'IPCC'
$data[6,8] = '[''region'', ''time_filter'']'
$data[7,0] = 11
$data[7,1] = 'future'
$data[7,2] = 'ts'
$data[7,3] = '[''ifs-fesom'', ''icon'']'
$data[7,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[7,5] = 'time_series/future_EERIE_EDDY_ts.zarr'
$data[7,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[7,7] = 'EDDY'
$data[7,8] = '[''region'', ''time_filter'']'
$data[8,0] = 1
$data[8,1] = 'hist'
$data[8,2] = 'clim'
$data[8,3] = '[''ifs-fesom'', ''icon'']'
$data[8,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[8,5] = 'decadal/hist_EERIE_clim.zarr'
$data[8,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[8,7] = ''
$data[8,8] = '[''period'', ''time_filter'']'
$data[9,0] = 6
$data[9,1] = 'hist'
$data[9,2] = 'trend'
$data[9,3] = '[''ifs-fesom'', ''icon'']'
$data[9,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[9,5] = 'decadal/hist_EERIE_trend.zarr'
$data[9,6] = '[''value (variable with hatching over variable_pvalue > 0.05)'']'
$data[9,7] = ''
$data[9,8] = '[''period'', ''time_filter'']'
$data[10,0] = 12
$data[10,1] = 'hist'
$data[10,2] = 'ts'
$data[10,3] = '[''ifs-fesom'', ''icon'']'
$data[10,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[10,5] = 'time_series/hist_EERIE_IPCC_ts.zarr'
$data[10,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[10,7] = 'IPCC'
$data[10,8] = '[''region'', ''time_filter'']'
$data[11,0] = 13
$data[11,1] = 'hist'
$data[11,2] = 'ts'
$data[11,3] = '[''ifs-fesom'', ''icon'']'
$data[11,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[11,5] = 'time_series/hist_EERIE_EDDY_ts.zarr'
$data[11,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[11,7] = 'EDDY'
$data[11,8] = '[''region'', ''time_filter'']'
$data[12,0] = 2
$data[12,1] = 'hist-amip'
$data[12,2] = 'clim'
$data[12,3] = '[''ifs-amip-tco1279.hist'', ''ifs-amip-tco1279.hist-c-0-a-lr20'', ''ifs-amip-tco399.hist-c-0-a-lr20'', ''ifs-amip-tco399.hist-c-lr20-a-0'', ''ifs-amip-tco399.hist'']'
$data[12,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''tasmax'', ''tasmin'']'
$data[12,5] = 'decadal/hist-amip_EERIE_clim.zarr'
$data[12,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[12,7] = ''
$data[12,8] = '[''period'', ''time_filter'']'
$data[13,0] = 7
$data[13,1] = 'hist-amip'
$data[13,2] = 'trend'
$data[13,3] = '[''ifs-amip-tco1279.hist'', ''ifs-amip-tco1279.hist-c-0-a-lr20'', ''ifs-amip-tco399.hist-c-0-a-lr20'', ''ifs-amip-tco399.hist-c-lr20-a-0'', ''ifs-amip-tco399.hist'']'
$data[13,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''tasmax'', ''tasmin'']'
$data[13,5] = 'decadal/hist-amip_EERIE_trend.zarr'
$data[13,6] = '[''value (variable with hatching over variable_pvalue > 0.05)'']'
$data[13,7] = ''
$data[13,8] = '[''period'', ''time_filter'']'
$data[14,0] = 14
$data[14,1] = 'hist-amip'
$data[14,2] = 'ts'
$data[14,3] = '[''ifs-amip-tco1279-hist'', ''ifs-amip-tco1279-hist-c-0-a-lr20'', ''ifs-amip-tco399-hist-c-0-a-lr20_1'', ''ifs-amip-tco399-hist-c-0-a-lr20_10'', ''ifs-amip-tco399-hist-c-0-a-lr20_2'', ''ifs-amip-tco399-hist-c-0-a-lr20_3'', ''ifs-amip-tco399-hist-c-0-a-lr20_4'', ''ifs-amip-tco399-hist-c-0-a-lr20_5'', ''ifs-amip-tco399-hist-c-0-a-lr20_6'', ''ifs-amip-tco399-hist-c-0-a-lr20_7'', ''ifs-amip-tco399-hist-c-0-a-lr20_8'', ''ifs-amip-tco399-hist-c-0-a-lr20_9'', ''ifs-amip-tco399-hist-c-lr20-a-0_1'', ''ifs-amip-tco399-hist-c-lr20-a-0_10'', ''ifs-amip-tco399-hist-c-lr20-a-0_2'', ''ifs-amip-tco399-hist-c-lr20-a-0_3'', ''ifs-amip-tco399-hist-c-lr20-a-0_4'', ''ifs-amip-tco399-hist-c-lr20-a-0_5'', ''ifs-amip-tco399-hist-c-lr20-a-0_6'', ''ifs-amip-tco399-hist-c-lr20-a-0_7'', ''ifs-amip-tco399-hist-c-lr20-a-0_8'', ''ifs-amip-tco399-hist-c-lr20-a-0_9'', ''ifs-amip-tco399-hist_1'', ''ifs-amip-tco399-hist_10'', ''ifs-amip-tco399-hist_2'', ''ifs-amip-tco399-hist_3'', ''ifs-amip-tco399-hist_4'', ''ifs-amip-tco399-hist_5'', ''ifs-amip-tco399-hist_6'', ''ifs-amip-tco399-hist_7'', ''ifs-amip-tco399-hist_8'', ''ifs-amip-tco399-hist_9'']'
$data[14,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''tasmax'', ''tasmin'']'
$data[14,5] = 'time_series/hist-amip_EERIE_IPCC_ts.zarr'
$data[14,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[14,7] = 'IPCC'
$data[14,8] = '[''region'', ''time_filter'']'
$data[15,0] = 15
$data[15,1] = 'hist-amip'
$data[15,2] = 'ts'
$data[15,3] = '[''ifs-amip-tco1279-hist'', ''ifs-amip-tco1279-hist-c-0-a-lr20'', ''ifs-amip-tco399-hist-c-0-a-lr20_1'', ''ifs-amip-tco399-hist-c-0-a-lr20_10'', ''ifs-amip-tco399-hist-c-0-a-lr20_2'', ''ifs-amip-tco399-hist-c-0-a-lr20_3'', ''ifs-amip-tco399-hist-c-0-a-lr20_4'', ''ifs-amip-tco399-hist-c-0-a-lr20_5'', ''ifs-amip-tco399-hist-c-0-a-lr20_6'', ''ifs-amip-tco399-hist-c-0-a-lr20_7'', ''ifs-amip-tco399-hist-c-0-a-lr20_8'', ''ifs-amip-tco399-hist-c-0-a-lr20_9'', ''ifs-amip-tco399-hist-c-lr20-a-0_1'', ''ifs-amip-tco399-hist-c-lr20-a-0_10'', ''ifs-amip-tco399-hist-c-lr20-a-0_2'', ''ifs-amip-tco399-hist-c-lr20-a-0_3'', ''ifs-amip-tco399-hist-c-lr20-a-0_4'', ''ifs-amip-tco399-hist-c-lr20-a-0_5'', ''ifs-amip-tco399-hist-c-lr20-a-0_6'', ''ifs-amip-tco399-hist-c-lr20-a-0_7'', ''ifs-amip-tco399-hist-c-lr20-a-0_8'', ''ifs-amip-tco399-hist-c-lr20-a-0_9'', ''ifs-amip-tco399-hist_1'', ''ifs-amip-tco399-hist_10'', ''ifs-amip-tco399-hist_2'', ''ifs-amip-tco399-hist_3'', ''ifs-amip-tco399-hist_4'', ''ifs-amip-tco399-hist_5'', ''ifs-amip-tco399-hist_6'', ''ifs-amip-tco399-hist_7'', ''ifs-amip-tco399-hist_8'', ''ifs-amip-tco399-hist_9'']'
$data[15,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''tasmax'', ''tasmin'']'
$data[15,5] = 'time_series/hist-amip_EERIE_EDDY_ts.zarr'
$data[15,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[15,7] = 'EDDY'
$data[15,8] = '[''region'', ''time_filter'']'
$data[16,0] = 4
$data[16,1] = 'obs'
$data[16,2] = 'clim'
$data[16,3] = '[''ifs-fesom'', ''icon'']'
$data[16,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[16,5] = 'decadal/obs_clim.zarr'
$data[16,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[16,7] = ''
$data[16,8] = '[''period'', ''time_filter'']'
$data[17,0] = 9
$data[17,1] = 'obs'
$data[17,2] = 'trend'
$data[17,3] = '[''ifs-fesom'', ''icon'']'
$data[17,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[17,5] = 'decadal/obs_trend.zarr'
$data[17,6] = '[''value (variable with hatching over variable_pvalue > 0.05)'']'
$data[17,7] = ''
$data[17,8] = '[''period'', ''time_filter'']'
$data[18,0] = 18
$data[18,1] = 'obs'
$data[18,2] = 'ts'
$data[18,3] = '[''ifs-fesom'', ''icon'']'
$data[18,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[18,5] = 'time_series/obs_IPCC_ts.zarr'
$data[18,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[18,7] = 'IPCC'
$data[18,8] = '[''region'', ''time_filter'']'
$data[19,0] = 19
$data[19,1] = 'obs'
$data[19,2] = 'ts'
$data[19,3] = '[''ifs-fesom'', ''icon'']'
$data[19,4] = '[''sfcWind'', ''uas'', ''vas'', ''tas'', ''pr'', ''tos'', ''clt'', ''zos'', ''tasmax'', ''tasmin'']'
$data[19,5] = 'time_series/obs_EDDY_ts.zarr'
$data[19,6] = '[''value (variable)'', ''anomaly (variable_anom)'']'
$data[19,7] = 'EDDY'
$data[19,8] = '[''region'', ''time_filter'']'

$ws.Range("A2:I21").Value = $data

# Apply style (bold, border, center/top alignment) to the new column-A cells
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A18:A21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0